$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(375, 44449, 5, 30, 90.95594700300154),
    @(376, 44450, 6, 33, 100.0515417033017),
    @(377, 44451, 3, 30, 90.95594700300154),
    @(378, 44452, 7, 27, 81.8603523027014),
    @(379, 44453, 0, 24, 72.76475760240123),
    @(380, 44454, 0, 23, 69.73289270230119),
    @(381, 44455, 7, 28, 84.89221720280145),
    @(382, 44456, 1, 24, 72.76475760240123),
    @(383, 44457, 5, 23, 69.73289270230119),
    @(384, 44458, 1, 21, 63.66916290210109),
    @(385, 44459, 2, 16, 48.50983840160082)
)

# Reference existing styled date cell (A374) so the new date cells (column A)
# inherit the same style (bordered, centered, date-formatted) used throughout
# the rest of the column.
$styleSource = $ws.Cells.Item(374, 1)

foreach ($row in $data) {
    $r = $row[0]

    # Copy formatting only (style/border/number-format/alignment) from A374,
    # then set the actual date value, so A375:A385 reuse the same cell style
    # as the rest of column A instead of creating a brand-new style entry.
    $cellA = $ws.Cells.Item($r, 1)
    $styleSource.Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null
    $cellA.Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
